$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B20 should become a real number (2) instead of a text string "2"
$ws.Range("B20").Value = 2

# Add new row 21 with the new annotation data
$ws.Range("A21").Value = "Ruilin"
# B21 holds the text "3" (not a number) - force text formatting so Excel
# doesn't auto-convert the numeric-looking string, then restore the
# default style so no extra formatting is left on the cell.
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = "3"
$ws.Range("B21").Style = "Normal"
$ws.Range("C21").Value = "无"
$ws.Range("D21").Value = "DFT"
$ws.Range("E21").Value = "WRI"
$ws.Range("F21").Value = "9386b51e-53f7-4a13-b66e-3217e88401e6"
$ws.Range("G21").Value = "HyRnez-RW_annotated.xlsx"
$ws.Range("H21").Value = '"Krasner" misspelled multiple times as "Kramer"'
